$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Detected Category" values in column C for several rows.
$ws.Range("C5").Value = "IMAGE_RETRIEVAL_BY_METADATA"
$ws.Range("C7").Value = "IMAGE_RETRIEVAL_BY_METADATA"

$ws.Range("C12").Value = "IMAGE_RETRIEVAL_BY_CAPTION"
$ws.Range("D12").Value = $true

$ws.Range("C14").Value = "IMAGE_RETRIEVAL_BY_METADATA"

$ws.Range("C21").Value = "IMAGE_RETRIEVAL_BY_METADATA"

$ws.Range("C25").Value = "IMAGE_RETRIEVAL_BY_METADATA"
